$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First initial balance run for levels 1-5
$ws.Range("B2").Value = 300
$ws.Range("C2").Value = 150

$ws.Range("C5").Value = 100
$ws.Range("C8").Value = 90
$ws.Range("C11").Value = 80
$ws.Range("C14").Value = 70
$ws.Range("C17").Value = 70
$ws.Range("C20").Value = 50

$ws.Range("C20").Select()
